# Added Molten Salt Storage.  Also fixed bug in battery calc
#
# The only real data change on the "Reference" sheet is the starting
# offset used for the per-energy-type parameter blocks (I2 / _1st_ref),
# which moves from 14 to 16 to make room for the new Molten Salt Storage
# parameters.  All of the K column values are computed via formulas
# referencing _1st_ref / Param_Count, so they recalculate automatically.
# The active selection on the sheet also moved from L1 to I3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reference")

# Update the "1st ref" value that _1st_ref points to (I2): 14 -> 16
$ws.Range("I2").Value = 16

# Update the active cell / selection shown in the sheet view: L1 -> I3
$ws.Range("I3").Select()

$wb.Application.Calculate()
